$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.776.66'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '2.405.33'
$ws.Range("E3").Value = '  -0.71%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''551.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.54%  '
$ws.Range("D6").Value = '''137.09'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  +3.67%  '
$ws.Range("E9").Value = '  -1.95%  '
$ws.Range("E10").Value = '  -2.03%  '
$ws.Range("E11").Value = '  -1.10%  '
$ws.Range("E12").Value = '  -1.95%  '
$ws.Range("D13").Value = '''25.35'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.74%  '
$ws.Range("D14").Value = '2.831.76'
$ws.Range("E14").Value = '  -0.73%  '
$ws.Range("D15").Value = '59.697.65'
$ws.Range("E15").Value = '  +0.07%  '
$ws.Range("E16").Value = '  -2.03%  '
$ws.Range("D17").Value = '2.390.95'
$ws.Range("E17").Value = '  -0.63%  '
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("E19").Value = '  -1.16%  '
$ws.Range("D20").Value = '''328.53'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.90%  '
$ws.Range("E21").Value = '  -4.17%  '
$ws.Range("D22").Value = '''0.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").Value = '''66.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.96%  '
$ws.Range("E24").Value = '  +2.06%  '
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("D27").Value = '''1.37'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.44%  '
$ws.Range("D28").Value = '0.0₃0771'
$ws.Range("E28").Value = '  -2.56%  '
$ws.Range("E29").Value = '  -2.23%  '
$ws.Range("D30").Value = '''168.37'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.18%  '
$ws.Range("D31").Value = '''6.03'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.18%  '
$ws.Range("D32").Value = '''18.60'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.70%  '
$ws.Range("E33").Value = '  -1.25%  '
$ws.Range("E35").Value = '  -0.79%  '
$ws.Range("E36").Value = '  -0.14%  '
$ws.Range("D37").Value = '''4.19'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.91%  '
$ws.Range("E38").Value = '  -2.05%  '
$ws.Range("D39").Value = '''319.76'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.03%  '
$ws.Range("D40").Value = '''0.407'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.22%  '
$ws.Range("D41").Value = '''3.67'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.20%  '
$ws.Range("D42").Value = '''139.31'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.35%  '
$ws.Range("D43").Value = '''0.0967'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.13%  '
$ws.Range("D44").Value = '''19.62'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.05%  '
$ws.Range("D45").Value = '''0.0514'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.86%  '
$ws.Range("E46").Value = '  +0.79%  '
$ws.Range("D47").Value = '''0.0223'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.52%  '
$ws.Range("D48").Value = '''0.386'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.25%  '
$ws.Range("D49").Value = '''17.54'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.43%  '
$ws.Range("D50").Value = '''11.04'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.03%  '
$ws.Range("D51").Value = '''1.57'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.04%  '
